$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2020" column (Q) mirroring the styles used in column P for each row.

# Row 4 - year header
$ws.Range("Q4").Value = 2020
$ws.Range("Q4").Style = $ws.Range("P4").Style

# Row 5
$ws.Range("Q5").Value = 99.3
$ws.Range("Q5").Style = $ws.Range("P5").Style

# Row 6
$ws.Range("Q6").Value = 99.371420589467803
$ws.Range("Q6").Style = $ws.Range("P6").Style

# Row 7
$ws.Range("Q7").Value = 99.319469393395053
$ws.Range("Q7").Style = $ws.Range("P7").Style

# Row 8
$ws.Range("Q8").Value = 99.442213297634979
$ws.Range("Q8").Style = $ws.Range("P8").Style

# Row 9
$ws.Range("Q9").Value = 98.766881972988841
$ws.Range("Q9").Style = $ws.Range("P9").Style

# Row 10
$ws.Range("Q10").Value = 99.212798374809537
$ws.Range("Q10").Style = $ws.Range("P10").Style

# Row 11
$ws.Range("Q11").Value = 99.799160124155549
$ws.Range("Q11").Style = $ws.Range("P11").Style

# Row 12
$ws.Range("Q12").Value = 99.146991622239156
$ws.Range("Q12").Style = $ws.Range("P12").Style

# Row 13
$ws.Range("Q13").Value = 99.538370126605429
$ws.Range("Q13").Style = $ws.Range("P13").Style

# Row 14
$ws.Range("Q14").Value = 99.765563948945029
$ws.Range("Q14").Style = $ws.Range("P14").Style

# Update the selection to match the authored state (cell P7 selected).
$ws.Range("P7").Select()
